$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted price strings (e.g. thousands-dot notation).
# Force text format before assigning so Excel does not auto-coerce numeric-looking
# strings (e.g. "586.86") into real numbers, which would change the stored cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '66.218.37'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '3.452.57'
$ws.Range('E3').Value = '  -1.12%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '586.86'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('D6').Value = '176.25'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.603'
$ws.Range('E8').Value = '  +0.70%  '
$ws.Range('D9').Value = '3.450.49'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('D10').Value = '0.132'
$ws.Range('E10').Value = '  -2.18%  '
$ws.Range('D11').Value = '6.98'
$ws.Range('E11').Value = '  +1.13%  '
$ws.Range('D12').Value = '0.418'
$ws.Range('E12').Value = '  -0.49%  '
$ws.Range('D13').Value = '4.052.29'
$ws.Range('E13').Value = '  -0.73%  '
$ws.Range('E14').Value = '  +1.57%  '
$ws.Range('D15').Value = '29.46'
$ws.Range('E15').Value = '  -3.36%  '
$ws.Range('D16').Value = '66.128.87'
$ws.Range('E16').Value = '  -0.82%  '
$ws.Range('D17').Value = '0.0000172'
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('D18').Value = '3.453.30'
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('D19').Value = '5.95'
$ws.Range('E19').Value = '  -1.35%  '
$ws.Range('D20').Value = '13.78'
$ws.Range('E20').Value = '  -1.42%  '
$ws.Range('D21').Value = '372.82'
$ws.Range('E21').Value = '  -1.65%  '
$ws.Range('D22').Value = '7.63'
$ws.Range('E22').Value = '  -2.10%  '
$ws.Range('D23').Value = '73.15'
$ws.Range('E23').Value = '  +2.40%  '
$ws.Range('D24').Value = '0.997'
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('D25').Value = '0.538'
$ws.Range('E25').Value = '  +1.33%  '
$ws.Range('E26').Value = '  +4.29%  '
$ws.Range('D27').Value = '9.78'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('D28').Value = '0.178'
$ws.Range('E28').Value = '  +2.87%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = '5.84'
$ws.Range('E30').Value = '  -1.65%  '
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('D32').Value = '23.70'
$ws.Range('E32').Value = '  -3.04%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('D34').Value = '7.05'
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('D35').Value = '1.27'
$ws.Range('E35').Value = '  -6.30%  '
$ws.Range('E36').Value = '  +1.89%  '
$ws.Range('D37').Value = '161.71'
$ws.Range('E37').Value = '  +1.80%  '
$ws.Range('E38').Value = '  +0.28%  '
$ws.Range('D39').Value = '28.38'
$ws.Range('E39').Value = '  +1.84%  '
$ws.Range('D40').Value = '1.81'
$ws.Range('E40').Value = '  +0.98%  '
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('D42').Value = '2.776.98'
$ws.Range('E42').Value = '  +2.66%  '
$ws.Range('D43').Value = '4.50'
$ws.Range('E43').Value = '  +0.34%  '
$ws.Range('E44').Value = '  -2.42%  '
$ws.Range('D45').Value = '0.0690'
$ws.Range('E45').Value = '  -1.36%  '
$ws.Range('D46').Value = '25.21'
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').Value = '336.67'
$ws.Range('E47').Value = '  +4.04%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = '39.90'
$ws.Range('E48').Value = '  -0.98%  '
$ws.Range('D49').Value = '0.0292'
$ws.Range('E49').Value = '  -1.13%  '
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').Value = '6.27'
$ws.Range('E51').Value = '  +1.17%  '
